$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "CDPDash"
$ws.Range("B3").Value = "N"

$ws.Range("B3").Select()
